{"js": "// Office.js (Word JavaScript API) script.\n// Applies two content edits described by the diff:\n//   1) In the time-log table, the task cell\n//      \"Implemented test cases to test bed main in ShoppingBag.java\"\n//      becomes\n//      \"Implement documented test cases to test bed main in ShoppingBag.java\"\n//      (the single run is split into \"Implement\" / \" \" / \"documented \" / \"test cases...\").\n//   2) The last paragraph of the document (\"The biggest issue with ...\") gets a\n//      reflection about comments/Javadoc inserted in the middle and a large\n//      new block of text appended at the end.\n//\n// NOTE: each `insertText(..., Word.InsertLocation.after)` call must chain off\n// the Range object *returned* by the previous insertText call (not the\n// original search-result range) - the host keeps \"after\" anchored to the\n// range it was called on, so re-using a stale range would insert everything\n// immediately after the original anchor instead of building up the text in\n// the intended left-to-right order.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Edit 1: table cell \"Implemented test cases to test bed main in ShoppingBag.java\"\n// ---------------------------------------------------------------------\nconst taskResults = body.search(\n  \"Implemented test cases to test bed main in ShoppingBag.java\",\n  { matchCase: true, matchWholeWord: false }\n);\ntaskResults.load(\"items\");\nawait context.sync();\n\nif (taskResults.items.length > 0) {\n  let taskRange = taskResults.items[0];\n  // Replace the whole run's text with just \"Implement\" (the new first run).\n  taskRange = taskRange.insertText(\"Implement\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Append the remaining pieces right after, in order, so the final visible\n  // text reads \"Implement documented test cases to test bed main in ShoppingBag.java\".\n  taskRange = taskRange.insertText(\" \", Word.InsertLocation.after);\n  await context.sync();\n  taskRange = taskRange.insertText(\"documented \", Word.InsertLocation.after);\n  await context.sync();\n  taskRange = taskRange.insertText(\n    \"test cases to test bed main in ShoppingBag.java\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Edit 2: closing reflection paragraph\n// ---------------------------------------------------------------------\n// 2a) Split \". For the project itself, formatting the proper outputs\" apart,\n//     inserting the new \"comment requirements\" sentence and turning\n//     \"For the project itself\" into \"For the project code itself\".\nconst midResults = body.search(\n  \". For the project itself, formatting the proper outputs\",\n  { matchCase: true }\n);\nmidResults.load(\"items\");\nawait context.sync();\n\nif (midResults.items.length > 0) {\n  let midRange = midResults.items[0];\n  midRange = midRange.insertText(\". \", Word.InsertLocation.replace);\n  await context.sync();\n\n  midRange = midRange.insertText(\n    \"While I had always commented code in the past, the new comment requirements took some time to adjust to. \",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n\n  midRange = midRange.insertText(\"For the project\", Word.InsertLocation.after);\n  await context.sync();\n\n  midRange = midRange.insertText(\" code\", Word.InsertLocation.after);\n  await context.sync();\n\n  midRange = midRange.insertText(\n    \" itself, formatting the proper outputs\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n\n// 2b) Append the new Javadoc / lessons-learned material at the very end of the\n//     closing sentence (\"... our own code.\").\nconst tailResults = body.search(\n  \"was time consuming as I had to repeated check between the project description and our own code.\",\n  { matchCase: true }\n);\ntailResults.load(\"items\");\nawait context.sync();\n\nif (tailResults.items.length > 0) {\n  let tailRange = tailResults.items[0];\n\n  tailRange = tailRange.insertText(\" \", Word.InsertLocation.after);\n  await context.sync();\n\n  tailRange = tailRange.insertText(\n    \"Outside of the project code, the biggest issue I had was generating the \",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n\n  tailRange = tailRange.insertText(\"Javadoc\", Word.InsertLocation.after);\n  await context.sync();\n\n  tailRange = tailRange.insertText(\n    \". This was due to a combination of \",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n\n  tailRange = tailRange.insertText(\n    \"my first attempt at creating a Java doc paired with a Java version issue that caused me to lose about 90 minutes of time trying to solve. The issue was resolved once a version of Java was installed that allowed the creation of \",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n\n  tailRange = tailRange.insertText(\"Javadocs\", Word.InsertLocation.after);\n  await context.sync();\n\n  tailRange = tailRange.insertText(\" as described during lecture. \", Word.InsertLocation.after);\n  await context.sync();\n\n  tailRange = tailRange.insertText(\n    \"While these issues took time to resolve, they were issue that rose from a lack of experience rather than lack of understanding and will no longer negatively affect me as much in the future. \",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies two content edits described by the diff:\n#   1) In the time-log table, the task cell\n#      \"Implemented test cases to test bed main in ShoppingBag.java\"\n#      becomes\n#      \"Implement documented test cases to test bed main in ShoppingBag.java\".\n#   2) The last paragraph of the document (\"The biggest issue with ...\") gets a\n#      reflection about comments/Javadoc inserted in the middle and a large\n#      new block of text appended at the end.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Edit 1: table cell \"Implemented test cases to test bed main in ShoppingBag.java\"\n# ---------------------------------------------------------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"Implemented test cases to test bed main in ShoppingBag.java\"\n$find1.Replacement.Text = \"Implement documented test cases to test bed main in ShoppingBag.java\"\n$find1.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# ---------------------------------------------------------------------\n# Edit 2: closing reflection paragraph\n# ---------------------------------------------------------------------\n# 2a) Split \". For the project itself, formatting the proper outputs\" apart,\n#     inserting the new \"comment requirements\" sentence and turning\n#     \"For the project itself\" into \"For the project code itself\".\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \". For the project itself, formatting the proper outputs\"\n$find2.Replacement.Text = \". While I had always commented code in the past, the new comment requirements took some time to adjust to. For the project code itself, formatting the proper outputs\"\n$find2.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# 2b) Append the new Javadoc / lessons-learned material at the very end of the\n#     closing sentence (\"... our own code.\").\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"was time consuming as I had to repeated check between the project description and our own code.\"\n$find3.Replacement.Text = \"was time consuming as I had to repeated check between the project description and our own code. Outside of the project code, the biggest issue I had was generating the Javadoc. This was due to a combination of my first attempt at creating a Java doc paired with a Java version issue that caused me to lose about 90 minutes of time trying to solve. The issue was resolved once a version of Java was installed that allowed the creation of Javadocs as described during lecture. While these issues took time to resolve, they were issue that rose from a lack of experience rather than lack of understanding and will no longer negatively affect me as much in the future. \"\n$find3.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2) | Out-Null\n"}
